$wb = $excel.ActiveWorkbook

# --- quality sheet (index 4) ---
$ws = $wb.Worksheets.Item(4)
$qualityData = @{
  2 = @{ "H" = "21,30316;31,10658;41,2567"; "I" = "21,0;31,31976;41,0"; "J" = "21,22644;31,7988;41,0"; "K" = "21,3683;31,1283;41,360"; "L" = "21,0;31,3850;41,0"; "M" = "21,2751;31,961;41,0"; "N" = 30316.4715; "O" = 10658.739; "P" = 2567.2174; "Q" = 0; "R" = 31976.217; "S" = 0; "T" = 22644.5481; "U" = 7988.1978; "V" = 0; "W" = 3683.673; "X" = 1283.373; "Y" = 360.3278; "Z" = 0; "AA" = 3850.119; "AB" = 0; "AC" = 2751.4782; "AD" = 961.8245999999999; "AE" = 0 }
  3 = @{ "H" = "21,30319;31,10660;41,2567"; "I" = "21,0;31,31981;41,0"; "J" = "21,22646;31,7989;41,0"; "K" = "21,3686;31,1285;41,360"; "L" = "21,0;31,3855;41,0"; "M" = "21,2753;31,963;41,0"; "N" = 30319.5585; "O" = 10660.377; "P" = 2567.6412; "Q" = 0; "R" = 31981.131; "S" = 0; "T" = 22646.8539; "U" = 7989.425399999999; "V" = 0; "W" = 3686.76; "X" = 1285.011; "Y" = 360.7515999999999; "Z" = 0; "AA" = 3855.033; "AB" = 0; "AC" = 2753.784; "AD" = 963.0522; "AE" = 0 }
  4 = @{ "H" = "21,32480;31,11816;41,2961"; "I" = "21,0;31,35448;41,0"; "J" = "21,24260;31,8855;41,0"; "K" = "21,3919;31,1409;41,404"; "L" = "21,0;31,4228;41,0"; "M" = "21,2927;31,1056;41,0"; "N" = 32480.5367236; "O" = 11816.087465; "P" = 2961.36720774; "Q" = 0; "R" = 35448.26239500001; "S" = 0; "T" = 24260.97232824; "U" = 8855.573243; "V" = 0; "W" = 3919.5236122; "X" = 1409.431205; "Y" = 404.24422822; "Z" = 0; "AA" = 4228.293615000001; "AB" = 0; "AC" = 2927.64416748; "AD" = 1056.298991; "AE" = 0 }
  5 = @{ "H" = "21,34933;31,13125;41,3407"; "I" = "21,0;31,39375;41,0"; "J" = "21,26093;31,9836;41,0"; "K" = "21,4183;31,1550;41,453"; "L" = "21,0;31,4650;41,0"; "M" = "21,3124;31,1161;41,0"; "N" = 34933.57737940001; "O" = 13125.0263599; "P" = 3407.4217477; "Q" = 0; "R" = 39375.0790797; "S" = 0; "T" = 26093.24351196; "U" = 9836.55821698; "V" = 0; "W" = 4183.3482313; "X" = 1550.1994963; "Y" = 453.5000181; "Z" = 0; "AA" = 4650.5984889; "AB" = 0; "AC" = 3124.70500542; "AD" = 1161.79786426; "AE" = 0 }
  6 = @{ "H" = "21,37816;31,14664;41,3932"; "I" = "21,0;31,43994;41,0"; "J" = "21,28246;31,10990;41,0"; "K" = "21,4493;31,1715;41,511"; "L" = "21,0;31,5147;41,0"; "M" = "21,3356;31,1285;41,0"; "N" = 37816.1499968; "O" = 14664.9253478; "P" = 3932.32455802; "Q" = 0; "R" = 43994.7760434; "S" = 0; "T" = 28246.34877312; "U" = 10990.63635956; "V" = 0; "W" = 4493.1925136; "X" = 1715.6698286; "Y" = 511.4249890599999; "Z" = 0; "AA" = 5147.0094858; "AB" = 0; "AC" = 3356.13971424; "AD" = 1285.80969572; "AE" = 0 }
  7 = @{ "H" = "21,42485;31,16438;41,5028"; "I" = "21,0;31,49314;41,0"; "J" = "21,31734;31,12319;41,0"; "K" = "21,5000;31,1908;41,633"; "L" = "21,0;31,5725;41,0"; "M" = "21,3734;31,1430;41,0"; "N" = 42485.7860175; "O" = 16438.3324287; "P" = 5028.89486236; "Q" = 0; "R" = 49314.9972861; "S" = 0; "T" = 31734.2809845; "U" = 12319.71727074; "V" = 0; "W" = 5000.12212875; "X" = 1908.3902019; "Y" = 633.2121570799999; "Z" = 0; "AA" = 5725.1706057; "AB" = 0; "AC" = 3734.78510025; "AD" = 1430.24408538; "AE" = 0 }
  8 = @{ "H" = "21,48526;31,18441;41,6545"; "I" = "21,0;31,55325;41,0"; "J" = "21,36246;31,13821;41,0"; "K" = "21,5653;31,2124;41,801"; "L" = "21,0;31,6374;41,0"; "M" = "21,4222;31,1592;41,0"; "N" = 48526.7608322; "O" = 18441.7896026; "P" = 6545.36526456; "Q" = 0; "R" = 55325.36880780001; "S" = 0; "T" = 36246.51931548; "U" = 13821.20935052; "V" = 0; "W" = 5653.2818069; "X" = 2124.9026162; "Y" = 801.27301368; "Z" = 0; "AA" = 6374.7078486; "AB" = 0; "AC" = 4222.65539046; "AD" = 1592.50943324; "AE" = 0 }
  9 = @{ "H" = "21,56293;31,20751;41,8586"; "I" = "21,0;31,62255;41,0"; "J" = "21,42047;31,15552;41,0"; "K" = "21,6489;31,2373;41,1026"; "L" = "21,0;31,7120;41,0"; "M" = "21,4847;31,1778;41,0"; "N" = 56293.28593000001; "O" = 20751.93383445; "P" = 8586.164629059998; "Q" = 0; "R" = 62255.80150335001; "S" = 0; "T" = 42047.63806200001; "U" = 15552.54821439; "V" = 0; "W" = 6489.952735000001; "X" = 2373.40386465; "Y" = 1026.90568218; "Z" = 0; "AA" = 7120.211593950001; "AB" = 0; "AC" = 4847.597349000001; "AD" = 1778.74883043; "AE" = 0 }
  10 = @{ "H" = "21,62113;31,21984;41,10285"; "I" = "21,0;31,65952;41,0"; "J" = "21,46395;31,16476;41,0"; "K" = "21,7116;31,2505;41,1214"; "L" = "21,0;31,7517;41,0"; "M" = "21,5315;31,1878;41,0"; "N" = 62113.69678500001; "O" = 21984.1447644; "P" = 10285.81436606; "Q" = 0; "R" = 65952.4342932; "S" = 0; "T" = 46395.128619; "U" = 16476.02937288; "V" = 0; "W" = 7116.9678825; "X" = 2505.8843628; "Y" = 1214.83604318; "Z" = 0; "AA" = 7517.6530884; "AB" = 0; "AC" = 5315.9392755; "AD" = 1878.03641256; "AE" = 0 }
  11 = @{ "H" = "21,69213;31,23447;41,12372"; "I" = "21,0;31,70341;41,0"; "J" = "21,51698;31,17572;41,0"; "K" = "21,7880;31,2662;41,1445"; "L" = "21,0;31,7988;41,0"; "M" = "21,5886;31,1995;41,0"; "N" = 69213.4435881; "O" = 23447.13378735; "P" = 12372.55046228; "Q" = 0; "R" = 70341.40136205; "S" = 0; "T" = 51698.20480254; "U" = 17572.46729997; "V" = 0; "W" = 7880.77192245; "X" = 2662.88490195; "Y" = 1445.35729684; "Z" = 0; "AA" = 7988.654705850001; "AB" = 0; "AC" = 5886.45412983; "AD" = 1995.70055289; "AE" = 0 }
  12 = @{ "H" = "21,82291;31,26526;41,16086"; "I" = "21,0;31,79579;41,0"; "J" = "21,61466;31,19880;41,0"; "K" = "21,9285;31,2993;41,1855"; "L" = "21,0;31,8979;41,0"; "M" = "21,6936;31,2243;41,0"; "N" = 82291.7296068; "O" = 26526.34026315; "P" = 16086.03012168; "Q" = 0; "R" = 79579.02078945; "S" = 0; "T" = 61466.88374712; "U" = 19880.18028513; "V" = 0; "W" = 9285.9023586; "X" = 2993.23406655; "Y" = 1855.12322504; "Z" = 0; "AA" = 8979.70219965; "AB" = 0; "AC" = 6936.00053724; "AD" = 2243.28091581; "AE" = 0 }
  13 = @{ "H" = "21,84316;31,27145;41,16612"; "I" = "21,0;31,81435;41,0"; "J" = "21,62979;31,20343;41,0"; "K" = "21,9510;31,3061;41,1914"; "L" = "21,0;31,9185;41,0"; "M" = "21,7103;31,2294;41,0"; "N" = 84316.4313432; "O" = 27145.1770772; "P" = 16612.334732; "Q" = 0; "R" = 81435.5312316; "S" = 0; "T" = 62979.21198288001; "U" = 20343.96787544; "V" = 0; "W" = 9510.2269164; "X" = 3061.7368964; "Y" = 1914.449996; "Z" = 0; "AA" = 9185.2106892; "AB" = 0; "AC" = 7103.55724776; "AD" = 2294.62039928; "AE" = 0 }
  14 = @{ "H" = "21,86341;31,27762;41,17138"; "I" = "21,0;31,83288;41,0"; "J" = "21,64491;31,20806;41,0"; "K" = "21,9734;31,3130;41,1973"; "L" = "21,0;31,9390;41,0"; "M" = "21,7271;31,2345;41,0"; "N" = 86341.1330796; "O" = 27762.9646931; "P" = 17138.4240943; "Q" = 0; "R" = 83288.89407930001; "S" = 0; "T" = 64491.54021864; "U" = 20806.96914362; "V" = 0; "W" = 9734.551474200001; "X" = 3130.1280647; "Y" = 1973.7822079; "Z" = 0; "AA" = 9390.384194100001; "AB" = 0; "AC" = 7271.113958280001; "AD" = 2345.87619794; "AE" = 0 }
  15 = @{ "H" = "21,88365;31,28381;41,17664"; "I" = "21,0;31,85145;41,0"; "J" = "21,66003;31,21270;41,0"; "K" = "21,9958;31,3198;41,2033"; "L" = "21,0;31,9596;41,0"; "M" = "21,7438;31,2397;41,0"; "N" = 88365.834816; "O" = 28381.84700715; "P" = 17664.76130462; "Q" = 0; "R" = 85145.54102145002; "S" = 0; "T" = 66003.8684544; "U" = 21270.79083393; "V" = 0; "W" = 9958.876032; "X" = 3198.67639455; "Y" = 2033.14157886; "Z" = 0; "AA" = 9596.029183650002; "AB" = 0; "AC" = 7438.670668799999; "AD" = 2397.24978141; "AE" = 0 }
  16 = @{ "H" = "21,90390;31,28999;41,18191"; "I" = "21,0;31,86998;41,0"; "J" = "21,67516;31,21733;41,0"; "K" = "21,10183;31,3267;41,2092"; "L" = "21,0;31,9801;41,0"; "M" = "21,7606;31,2448;41,0"; "N" = 90390.53655240001; "O" = 28999.63462305; "P" = 18191.09851494; "Q" = 0; "R" = 86998.90386915; "S" = 0; "T" = 67516.19669016; "U" = 21733.79210211; "V" = 0; "W" = 10183.2005898; "X" = 3267.06756285; "Y" = 2092.50094982; "Z" = 0; "AA" = 9801.202688550002; "AB" = 0; "AC" = 7606.22737932; "AD" = 2448.50558007; "AE" = 0 }
  17 = @{ "H" = "21,92418;31,29618;41,18717"; "I" = "21,0;31,88855;41,0"; "J" = "21,69030;31,22197;41,0"; "K" = "21,10407;31,3335;41,2151"; "L" = "21,0;31,10006;41,0"; "M" = "21,7774;31,2499;41,0"; "N" = 92418.1942677; "O" = 29618.4714371; "P" = 18717.40312526; "Q" = 0; "R" = 88855.4143113; "S" = 0; "T" = 69030.73286118; "U" = 22197.57969242; "V" = 0; "W" = 10407.81784665; "X" = 3335.5703927; "Y" = 2151.82772078; "Z" = 0; "AA" = 10006.7111781; "AB" = 0; "AC" = 7774.002718109999; "AD" = 2499.84506354; "AE" = 0 }
}
foreach ($r in $qualityData.Keys) {
  $rowdata = $qualityData[$r]
  foreach ($col in $rowdata.Keys) {
    $addr = "$col$r"
    $ws.Range($addr).Value = $rowdata[$col]
  }
}

# --- limiter sheet (index 7) ---
$ws = $wb.Worksheets.Item(7)
$limiterData = @{
  2 = @{ "B" = 5175; "C" = 1658; "D" = 1050; "E" = 28; "F" = 52; "G" = 37; "H" = 0; "I" = 3300; "J" = 0; "K" = 0; "L" = 156; "M" = 0; "N" = 3555; "O" = 1139; "P" = 0; "Q" = 21; "R" = 39; "S" = 0; "T" = 2633; "U" = 924; "V" = 417; "W" = 24; "X" = 45; "Y" = 32; "Z" = 0; "AA" = 1839; "AB" = 0; "AC" = 0; "AD" = 136; "AE" = 0; "AF" = 1809; "AG" = 634; "AH" = 0; "AI" = 17; "AJ" = 33; "AK" = 0; "AL" = "21,5175;31,1658;41,1050;22,28;32,52;42,37"; "AM" = "21,0;31,3300;41,0;22,0;32,156;42,0"; "AN" = "21,3555;31,1139;41,0;22,21;32,39;42,0"; "AO" = "21,2633;31,924;41,417;22,24;32,45;42,32"; "AP" = "21,0;31,1839;41,0;22,0;32,136;42,0"; "AQ" = "21,1809;31,634;41,0;22,17;32,33;42,0" }
  3 = @{ "B" = 10350; "C" = 3317; "D" = 2100; "E" = 56; "F" = 104; "G" = 74; "H" = 0; "I" = 6601; "J" = 0; "K" = 0; "L" = 312; "M" = 0; "N" = 7110; "O" = 2278; "P" = 0; "Q" = 42; "R" = 78; "S" = 0; "T" = 3040; "U" = 1026; "V" = 520; "W" = 49; "X" = 91; "Y" = 65; "Z" = 0; "AA" = 2041; "AB" = 0; "AC" = 0; "AD" = 273; "AE" = 0; "AF" = 2089; "AG" = 704; "AH" = 0; "AI" = 36; "AJ" = 68; "AK" = 0; "AL" = "21,10350;31,3317;41,2100;22,56;32,104;42,74"; "AM" = "21,0;31,6601;41,0;22,0;32,312;42,0"; "AN" = "21,7110;31,2278;41,0;22,42;32,78;42,0"; "AO" = "21,3040;31,1026;41,520;22,49;32,91;42,65"; "AP" = "21,0;31,2041;41,0;22,0;32,273;42,0"; "AQ" = "21,2089;31,704;41,0;22,36;32,68;42,0" }
  4 = @{ "B" = 15526; "C" = 4975; "D" = 3151; "E" = 84; "F" = 156; "G" = 112; "H" = 0; "I" = 9902; "J" = 0; "K" = 0; "L" = 468; "M" = 0; "N" = 10666; "O" = 3418; "P" = 0; "Q" = 63; "R" = 117; "S" = 0; "T" = 3734; "U" = 1203; "V" = 698; "W" = 73; "X" = 136; "Y" = 98; "Z" = 0; "AA" = 2395; "AB" = 0; "AC" = 0; "AD" = 409; "AE" = 0; "AF" = 2565; "AG" = 826; "AH" = 0; "AI" = 54; "AJ" = 101; "AK" = 0; "AL" = "21,15526;31,4975;41,3151;22,84;32,156;42,112"; "AM" = "21,0;31,9902;41,0;22,0;32,468;42,0"; "AN" = "21,10666;31,3418;41,0;22,63;32,117;42,0"; "AO" = "21,3734;31,1203;41,698;22,73;32,136;42,98"; "AP" = "21,0;31,2395;41,0;22,0;32,409;42,0"; "AQ" = "21,2565;31,826;41,0;22,54;32,101;42,0" }
  5 = @{ "B" = 20701; "C" = 6634; "D" = 4201; "E" = 112; "F" = 208; "G" = 149; "H" = 0; "I" = 13202; "J" = 0; "K" = 0; "L" = 624; "M" = 0; "N" = 14221; "O" = 4557; "P" = 0; "Q" = 84; "R" = 156; "S" = 0; "T" = 4005; "U" = 1289; "V" = 761; "W" = 98; "X" = 182; "Y" = 130; "Z" = 0; "AA" = 2566; "AB" = 0; "AC" = 0; "AD" = 546; "AE" = 0; "AF" = 2751; "AG" = 885; "AH" = 0; "AI" = 73; "AJ" = 136; "AK" = 0; "AL" = "21,20701;31,6634;41,4201;22,112;32,208;42,149"; "AM" = "21,0;31,13202;41,0;22,0;32,624;42,0"; "AN" = "21,14221;31,4557;41,0;22,84;32,156;42,0"; "AO" = "21,4005;31,1289;41,761;22,98;32,182;42,130"; "AP" = "21,0;31,2566;41,0;22,0;32,546;42,0"; "AQ" = "21,2751;31,885;41,0;22,73;32,136;42,0" }
  6 = @{ "B" = 25877; "C" = 8293; "D" = 5252; "E" = 140; "F" = 260; "G" = 186; "H" = 0; "I" = 16503; "J" = 0; "K" = 0; "L" = 780; "M" = 0; "N" = 17777; "O" = 5697; "P" = 0; "Q" = 105; "R" = 195; "S" = 0; "T" = 4293; "U" = 1380; "V" = 828; "W" = 122; "X" = 227; "Y" = 163; "Z" = 0; "AA" = 2747; "AB" = 0; "AC" = 0; "AD" = 682; "AE" = 0; "AF" = 2949; "AG" = 948; "AH" = 0; "AI" = 91; "AJ" = 169; "AK" = 0; "AL" = "21,25877;31,8293;41,5252;22,140;32,260;42,186"; "AM" = "21,0;31,16503;41,0;22,0;32,780;42,0"; "AN" = "21,17777;31,5697;41,0;22,105;32,195;42,0"; "AO" = "21,4293;31,1380;41,828;22,122;32,227;42,163"; "AP" = "21,0;31,2747;41,0;22,0;32,682;42,0"; "AQ" = "21,2949;31,948;41,0;22,91;32,169;42,0" }
  7 = @{ "B" = 31052; "C" = 9951; "D" = 6302; "E" = 168; "F" = 312; "G" = 224; "H" = 0; "I" = 19804; "J" = 0; "K" = 0; "L" = 936; "M" = 0; "N" = 21332; "O" = 6836; "P" = 0; "Q" = 126; "R" = 234; "S" = 0; "T" = 4601; "U" = 1477; "V" = 900; "W" = 147; "X" = 273; "Y" = 196; "Z" = 0; "AA" = 2940; "AB" = 0; "AC" = 0; "AD" = 819; "AE" = 0; "AF" = 3160; "AG" = 1015; "AH" = 0; "AI" = 110; "AJ" = 204; "AK" = 0; "AL" = "21,31052;31,9951;41,6302;22,168;32,312;42,224"; "AM" = "21,0;31,19804;41,0;22,0;32,936;42,0"; "AN" = "21,21332;31,6836;41,0;22,126;32,234;42,0"; "AO" = "21,4601;31,1477;41,900;22,147;32,273;42,196"; "AP" = "21,0;31,2940;41,0;22,0;32,819;42,0"; "AQ" = "21,3160;31,1015;41,0;22,110;32,204;42,0" }
  8 = @{ "B" = 36227; "C" = 11610; "D" = 7352; "E" = 196; "F" = 364; "G" = 261; "H" = 0; "I" = 23104; "J" = 0; "K" = 0; "L" = 1092; "M" = 0; "N" = 24887; "O" = 7975; "P" = 0; "Q" = 147; "R" = 273; "S" = 0; "T" = 4928; "U" = 1581; "V" = 977; "W" = 171; "X" = 318; "Y" = 228; "Z" = 0; "AA" = 3146; "AB" = 0; "AC" = 0; "AD" = 955; "AE" = 0; "AF" = 3385; "AG" = 1086; "AH" = 0; "AI" = 128; "AJ" = 237; "AK" = 0; "AL" = "21,36227;31,11610;41,7352;22,196;32,364;42,261"; "AM" = "21,0;31,23104;41,0;22,0;32,1092;42,0"; "AN" = "21,24887;31,7975;41,0;22,147;32,273;42,0"; "AO" = "21,4928;31,1581;41,977;22,171;32,318;42,228"; "AP" = "21,0;31,3146;41,0;22,0;32,955;42,0"; "AQ" = "21,3385;31,1086;41,0;22,128;32,237;42,0" }
  9 = @{ "B" = 41403; "C" = 13268; "D" = 8403; "E" = 220; "F" = 409; "G" = 294; "H" = 0; "I" = 26405; "J" = 0; "K" = 0; "L" = 1228; "M" = 0; "N" = 28443; "O" = 9115; "P" = 0; "Q" = 164; "R" = 306; "S" = 0; "T" = 5476; "U" = 1755; "V" = 1108; "W" = 196; "X" = 364; "Y" = 261; "Z" = 0; "AA" = 3493; "AB" = 0; "AC" = 0; "AD" = 1092; "AE" = 0; "AF" = 3762; "AG" = 1205; "AH" = 0; "AI" = 147; "AJ" = 273; "AK" = 0; "AL" = "21,41403;31,13268;41,8403;22,220;32,409;42,294"; "AM" = "21,0;31,26405;41,0;22,0;32,1228;42,0"; "AN" = "21,28443;31,9115;41,0;22,164;32,306;42,0"; "AO" = "21,5476;31,1755;41,1108;22,196;32,364;42,261"; "AP" = "21,0;31,3493;41,0;22,0;32,1092;42,0"; "AQ" = "21,3762;31,1205;41,0;22,147;32,273;42,0" }
  10 = @{ "B" = 46578; "C" = 14927; "D" = 9453; "E" = 234; "F" = 435; "G" = 312; "H" = 0; "I" = 29706; "J" = 0; "K" = 0; "L" = 1306; "M" = 0; "N" = 31998; "O" = 10254; "P" = 0; "Q" = 175; "R" = 325; "S" = 0; "T" = 5930; "U" = 1900; "V" = 1214; "W" = 220; "X" = 409; "Y" = 294; "Z" = 0; "AA" = 3782; "AB" = 0; "AC" = 0; "AD" = 1228; "AE" = 0; "AF" = 4073; "AG" = 1305; "AH" = 0; "AI" = 164; "AJ" = 306; "AK" = 0; "AL" = "21,46578;31,14927;41,9453;22,234;32,435;42,312"; "AM" = "21,0;31,29706;41,0;22,0;32,1306;42,0"; "AN" = "21,31998;31,10254;41,0;22,175;32,325;42,0"; "AO" = "21,5930;31,1900;41,1214;22,220;32,409;42,294"; "AP" = "21,0;31,3782;41,0;22,0;32,1228;42,0"; "AQ" = "21,4073;31,1305;41,0;22,164;32,306;42,0" }
  11 = @{ "B" = 51754; "C" = 16586; "D" = 10504; "E" = 245; "F" = 455; "G" = 326; "H" = 0; "I" = 33007; "J" = 0; "K" = 0; "L" = 1365; "M" = 0; "N" = 35554; "O" = 11394; "P" = 0; "Q" = 183; "R" = 341; "S" = 0; "T" = 6393; "U" = 2048; "V" = 1324; "W" = 245; "X" = 455; "Y" = 326; "Z" = 0; "AA" = 4077; "AB" = 0; "AC" = 0; "AD" = 1365; "AE" = 0; "AF" = 4392; "AG" = 1407; "AH" = 0; "AI" = 183; "AJ" = 341; "AK" = 0; "AL" = "21,51754;31,16586;41,10504;22,245;32,455;42,326"; "AM" = "21,0;31,33007;41,0;22,0;32,1365;42,0"; "AN" = "21,35554;31,11394;41,0;22,183;32,341;42,0"; "AO" = "21,6393;31,2048;41,1324;22,245;32,455;42,326"; "AP" = "21,0;31,4077;41,0;22,0;32,1365;42,0"; "AQ" = "21,4392;31,1407;41,0;22,183;32,341;42,0" }
}
foreach ($r in $limiterData.Keys) {
  $rowdata = $limiterData[$r]
  foreach ($col in $rowdata.Keys) {
    $addr = "$col$r"
    $ws.Range($addr).Value = $rowdata[$col]
  }
}

# --- status sheet (index 8) ---
$ws = $wb.Worksheets.Item(8)
$statusData = @{
  2 = @{ "T" = 5377713.924030701; "U" = 940849.9771480999; "V" = 818784.4610907601; "AJ" = "21,0;31,121862;41,0"; "AK" = "21,104584;31,33591;41,0"; "AR" = 631006.115153935; "AS" = 16592979.41762446 }
  3 = @{ "T" = 677027.08776015; "U" = 134481.7418997; "V" = 107426.69574228; "AJ" = "21,0;31,14083;41,0"; "AK" = "21,12166;31,3906;41,0"; "AR" = 73933.46295076501; "AS" = 2238044.60738619 }
  7 = @{ "T" = 61987.810056; "U" = 12049.54768; "V" = 6816.1587258; "AJ" = "21,0;31,3417;41,0"; "AK" = "21,3638;31,1168;41,0"; "AR" = 23127.32622135; "AS" = 189478.10172135 }
  8 = @{ "T" = 67474.225212; "U" = 13785.33736; "V" = 7904.0736216; "AJ" = "21,0;31,6835;41,0"; "AK" = "21,7277;31,2336;41,0"; "AR" = 46261.55244270001; "AS" = 213196.9763202 }
  10 = @{ "T" = 247984.08419792; "U" = 43691.9920992; "V" = 29541.98565556; "AJ" = "21,0;31,3794;41,0"; "AK" = "21,3938;31,1262;41,0"; "AR" = 24795.25906737; "AS" = 731643.24720187 }
  11 = @{ "T" = 279862.2853816; "U" = 52571.04737599999; "V" = 35250.91035968; "AJ" = "21,0;31,18975;41,0"; "AK" = "21,19694;31,6315;41,0"; "AR" = 124023.84636466; "AS" = 857783.2468441599 }
  12 = @{ "T" = 288993.1687016; "U" = 54155.832176; "V" = 36411.86929968; "AJ" = "21,0;31,18975;41,0"; "AK" = "21,19694;31,6315;41,0"; "AR" = 124023.84636466; "AS" = 884694.65918916 }
  13 = @{ "T" = 325918.5375188; "U" = 64651.683048; "V" = 43149.89674288; "AJ" = "21,0;31,37333;41,0"; "AK" = "21,38895;31,12474;41,0"; "AR" = 245372.83165206; "AS" = 1032955.05682156 }
  14 = @{ "T" = 298124.0520216; "U" = 55740.616976; "V" = 37572.82823968; "AJ" = "21,0;31,18975;41,0"; "AK" = "21,19694;31,6315;41,0"; "AR" = 124023.84636466; "AS" = 911423.0715341601 }
  15 = @{ "T" = 441274.69383888; "U" = 81177.68453832; "V" = 61842.48291180001; "AJ" = "21,0;31,24244;41,0"; "AK" = "21,24784;31,7946;41,0"; "AR" = 155792.176426688; "AS" = 1370005.993896138 }
  16 = @{ "T" = 447216.90218736; "U" = 83083.37848504001; "V" = 63036.8460227; "AJ" = "21,0;31,28284;41,0"; "AK" = "21,28914;31,9270;41,0"; "AR" = 181751.260895211; "AS" = 1396148.078364661 }
  17 = @{ "T" = 1077763.2514052; "U" = 191372.00982796; "V" = 170030.48084; "AJ" = "21,0;31,39304;41,0"; "AK" = "21,38415;31,12329;41,0"; "AR" = 240899.906249324; "AS" = 3389831.384048124 }
  18 = @{ "T" = 5377713.924030701; "U" = 940849.9771480999; "V" = 818784.4610907601; "AJ" = "21,0;31,121862;41,0"; "AK" = "21,104584;31,33591;41,0"; "AR" = 631006.115153935; "AS" = 16592979.41762446 }
  19 = @{ "T" = 8795360.848459201; "U" = 1523415.8256061; "V" = 1326611.76067126; "AJ" = "21,0;31,121862;41,0"; "AK" = "21,104584;31,33591;41,0"; "AR" = 631006.115153935; "AS" = 26952716.66900104 }
  20 = @{ "T" = 220081.607247; "U" = 41643.8132196; "V" = 25563.75339768; "AJ" = "21,0;31,2322;41,0"; "AK" = "21,2128;31,751;41,0"; "AR" = 13373.2620078; "AS" = 659571.0004989 }
  21 = @{ "T" = 286331.66922374; "U" = 52143.77200206; "V" = 35529.0148978; "AJ" = "21,0;31,3146;41,0"; "AK" = "21,2891;31,979;41,0"; "AR" = 17822.527268104; "AS" = 857997.781700304 }
  22 = @{ "T" = 399257.57681312; "U" = 70702.33664807; "V" = 53260.07458068; "AJ" = "21,0;31,6078;41,0"; "AK" = "21,5509;31,1773;41,0"; "AR" = 33420.118749263; "AS" = 1206007.128523713 }
  23 = @{ "T" = 459216.80755772; "U" = 83012.31128303999; "V" = 63845.9041797; "AJ" = "21,0;31,8205;41,0"; "AK" = "21,7301;31,2341;41,0"; "AR" = 44159.341590246; "AS" = 1412309.704443971 }
  24 = @{ "T" = 530612.152824; "U" = 99873.65863291999; "V" = 78209.88155353; "AJ" = "21,0;31,9995;41,0"; "AK" = "21,8828;31,2830;41,0"; "AR" = 53414.4143494705; "AS" = 1683373.216323945 }
  25 = @{ "T" = 677027.08776015; "U" = 134481.7418997; "V" = 107426.69574228; "AJ" = "21,0;31,14083;41,0"; "AK" = "21,12166;31,3906;41,0"; "AR" = 73933.46295076501; "AS" = 2238044.60738619 }
}
foreach ($r in $statusData.Keys) {
  $rowdata = $statusData[$r]
  foreach ($col in $rowdata.Keys) {
    $addr = "$col$r"
    $ws.Range($addr).Value = $rowdata[$col]
  }
}